$wb = $excel.ActiveWorkbook

# Source sheet to borrow cell formatting from (getTestInfoGridLayout)
$srcWs = $wb.Worksheets.Item(1)
# Insert new sheet after putDonationInfogridLayout (the last sheet)
$afterWs = $wb.Worksheets.Item(3)

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterWs)
$ws.Name = "putTestInfogridLayout"

# Copy the cell formatting (styles) from the getTestInfoGridLayout sheet,
# cell by cell so each cell keeps its own distinct style.
$cellsToFormat = @(
    "A1","B1","C1","D1","E1",
    "A2","B2","C2","D2","E2",
    "A3","B3","C3","D3","E3",
    "A4","B4","C4","D4","E4",
    "A5","B5","C5","D5","E5",
    "A6","B6","C6","D6","E6",
    "A7","B7","C7","D7","E7",
    "A8","B8","C8","D8","E8",
    "A9","B9","C9","D9","E9",
    "A10","B10","C10","D10","E10",
    "A11","B11","C11","D11","E11"
)
foreach ($addr in $cellsToFormat) {
    $srcWs.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# Values
$ws.Range("A1").Value = "Assert200"

$ws.Range("A2").Value = "gridName"
$ws.Range("B2").Value = "Body"
$ws.Range("C2").Value = "Page"
$ws.Range("D2").Value = "Sort"
$ws.Range("E2").Value = "EndPoint"

$ws.Range("A3").Value = "TestInformation"
$ws.Range("E3").Value = "/gridLayout/columnPreference"

$ws.Range("A5").Value = "Assert400"

$ws.Range("A6").Value = "gridName"
$ws.Range("B6").Value = "PageSize"
$ws.Range("C6").Value = "Page"
$ws.Range("D6").Value = "Sort"
$ws.Range("E6").Value = "EndPoint"

$ws.Range("A7").Value = "TestInformation"
$ws.Range("E7").Value = "/gridLayout/columnPreference"

$ws.Range("A9").Value = "Assert401"

$ws.Range("A10").Value = "gridName"
$ws.Range("B10").Value = "PageSize"
$ws.Range("C10").Value = "Page"
$ws.Range("D10").Value = "Sort"
$ws.Range("E10").Value = "EndPoint"

$ws.Range("A11").Value = "TestInformation"
$ws.Range("E11").Value = "/gridLayout/columnPreference"

# Row heights
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 50.1
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 50.1

# Column width (column E)
$ws.Columns.Item(5).ColumnWidth = 31

# Merge cells for the section header rows
$ws.Range("A1:E1").Merge()
$ws.Range("A5:E5").Merge()
$ws.Range("A9:E9").Merge()

# Sheet view settings
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("A11").Select()
